$wb = $excel.ActiveWorkbook

# --- Update Hoja1!A1 conversion summary text ---
$ws1 = $wb.Worksheets.Item("Hoja1")

$nl = [char]10
$newText = "Conversión del día 💰" + $nl + `
    "✅ Dólar paralelo: 68" + $nl + `
    $nl + `
    "Binance" + $nl + `
    "✅ 1000 Bs = 12.96 = 52986.14 pesos" + $nl + `
    "✅ 52986.14 pesos = 12.97 = 980.73 Bs" + $nl + `
    $nl + `
    "Promedio competencia" + $nl + `
    "✅ Tasa pesos: 20" + $nl + `
    "✅ Tasa Bs: 20" + $nl + `
    "✅ % Ganancia: 20%"

$ws1.Range("A1").Value = $newText

# --- Update tasas sheet rate values ---
$ws2 = $wb.Worksheets.Item("tasas")

$ws2.Range("N10").Value = 77.19
$ws2.Range("O10").Value = 4090
$ws2.Range("N12").Value = 4085
$ws2.Range("O12").Value = 75.61
